# Insert two new feature-comparison rows into the Top20 concat table.
#
# 1) "franzosa_ControlvsCD_ConvCD"  -> inserted right after the existing
#    "franzosa_ControlvsCD_Age" row (originally row 7), pushing the old
#    "franzosa_ControlvsCD_Fp" row (and everything below it) down by one.
# 2) "franzosa_ControlvsUC_ConvUC"  -> inserted right after the (now
#    shifted) "franzosa_ControlvsUC_Age" row, pushing the old
#    "franzosa_ControlvsUC_Fp" row (and everything below it) down by one
#    more.
#
# Net effect: the sheet grows from 24 data+header rows (A1:H24) to 26
# rows (A1:H26), and all rows from the old "franzosa_ControlvsDisease_Age"
# row through "wang_urea" end up two rows lower than before.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Insert row for franzosa_ControlvsCD_ConvCD (goes in at row 8) -------
$ws.Rows.Item(8).Insert()
$ws.Range("A8").Value = "franzosa_ControlvsCD_ConvCD"
$ws.Range("B8").Value = 0
$ws.Range("C8").Value = 0.05
$ws.Range("D8").Value = 0
$ws.Range("E8").Value = 0.45
$ws.Range("F8").Value = 0.95
$ws.Range("G8").Value = 0.5
$ws.Range("H8").Value = 0.55

# --- Insert row for franzosa_ControlvsUC_ConvUC (goes in at row 14) ------
$ws.Rows.Item(14).Insert()
$ws.Range("A14").Value = "franzosa_ControlvsUC_ConvUC"
$ws.Range("B14").Value = 0
$ws.Range("C14").Value = 0
$ws.Range("D14").Value = 0
$ws.Range("E14").Value = 0.35
$ws.Range("F14").Value = 1
$ws.Range("G14").Value = 0.65
$ws.Range("H14").Value = 0.65
